# ChkPoint: clean up scripts a4.6_all_data_preparation.py a5.2_pd_model_data_prep.py
#
# This script applies the edits described by the diff:
#   - Rename worksheet "independent_variables" to "PD_variables"
#   - Insert a new header row at the top of the sheet with headers
#     "COLUMN" (A1) and "Is Ref" (B1)
#   - Populate the new column B with indicator values (0 for every
#     dummy-variable row except the reference category "grade:G",
#     which is now on row 8 and gets the value 1)
#   - Update the active selection to B9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first worksheet.
$ws.Name = "PD_variables"

# Insert a new row 1, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# New header row.
$ws.Range("A1").Value = "COLUMN"
$ws.Range("B1").Value = "Is Ref"

# Fill in the "Is Ref" flag column for each of the data rows (rows 2-11).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}

# "grade:G" (now row 8) is the reference category -> flag it with 1.
$ws.Cells.Item(8, 2).Value = 1

# Match the saved selection/active cell.
$ws.Range("B9").Select()
